# Update "想去人数" (want-to-go count) figures that were re-scraped.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 1765
$wsExpo.Range("F15").Value = 663
$wsExpo.Range("F21").Value = 129
$wsExpo.Range("F23").Value = 25
$wsExpo.Range("F29").Value = 147

# --- Sheet "演出" (Shows) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 12

# --- Sheet "全部类型" (All types, merged view) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 1765
$wsAll.Range("F17").Value = 663
$wsAll.Range("F24").Value = 12
$wsAll.Range("F29").Value = 129
$wsAll.Range("F31").Value = 25
$wsAll.Range("F33").Value = 0
$wsAll.Range("F39").Value = 147

# Row 44 (the duplicate "KANAKO ITO&AYANE" event) was removed from the feed.
# Shift the content of B45:I49 up into B44:I48 (column A holds a positional
# index and is intentionally left untouched), then delete the now-duplicated
# trailing row 49 so the sheet shrinks from 49 to 48 data rows.
$wsAll.Range("B45:I49").Copy()
$wsAll.Range("B44").PasteSpecial()
$excel.CutCopyMode = $false
$wsAll.Rows.Item(49).Delete()
